$d = $word.ActiveDocument
$r = $d.Content
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se wp14"><w:body>
    <w:p w:rsidR="00763511" w:rsidRDefault="002F3A52" w:rsidP="002F3A52">
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:r>
        <w:t>Hola mundo</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> haremos un </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>commint</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r w:rsidR="00B21F90">
        <w:rPr>
          <w:noProof/>
          <w:lang w:eastAsia="es-HN"/>
        </w:rPr>
        <w:drawing>
          <wp:inline distT="0" distB="0" distL="0" distR="0">
            <wp:extent cx="5486400" cy="3200400"/>
            <wp:effectExtent l="0" t="0" r="19050" b="0"/>
            <wp:docPr id="1" name="Diagrama 1"/>
            <wp:cNvGraphicFramePr/>
            <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
              <a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/diagram">
                <dgm:relIds xmlns:dgm="http://schemas.openxmlformats.org/drawingml/2006/diagram" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" r:dm="rId4" r:lo="rId5" r:qs="rId6" r:cs="rId7"/>
              </a:graphicData>
            </a:graphic>
          </wp:inline>
        </w:drawing>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="left"/>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:rPr>
          <w:noProof/>
        </w:rPr>
        <mc:AlternateContent>
          <mc:Choice Requires="wps">
            <w:drawing>
              <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="58297744" wp14:editId="49D15B7F">
                <wp:simplePos x="0" y="0"/>
                <wp:positionH relativeFrom="column">
                  <wp:posOffset>0</wp:posOffset>
                </wp:positionH>
                <wp:positionV relativeFrom="paragraph">
                  <wp:posOffset>0</wp:posOffset>
                </wp:positionV>
                <wp:extent cx="1828800" cy="1828800"/>
                <wp:effectExtent l="0" t="0" r="0" b="0"/>
                <wp:wrapNone/>
                <wp:docPr id="2" name="Cuadro de texto 2"/>
                <wp:cNvGraphicFramePr/>
                <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
                  <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
                    <wps:wsp>
                      <wps:cNvSpPr txBox="1"/>
                      <wps:spPr>
                        <a:xfrm>
                          <a:off x="0" y="0"/>
                          <a:ext cx="1828800" cy="1828800"/>
                        </a:xfrm>
                        <a:prstGeom prst="rect">
                          <a:avLst/>
                        </a:prstGeom>
                        <a:noFill/>
                        <a:ln>
                          <a:noFill/>
                        </a:ln>
                      </wps:spPr>
                      <wps:txbx>
                        <w:txbxContent>
                          <w:p>
                            <w:pPr>
                              <w:jc w:val="center"/>
                              <w:rPr>
                                <w:b/>
                                <w:color w:val="F7CAAC" w:themeColor="accent2" w:themeTint="66"/>
                                <w:sz w:val="72"/>
                                <w:szCs w:val="72"/>
                                <w14:textOutline w14:w="11112" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
                                  <w14:solidFill>
                                    <w14:schemeClr w14:val="accent2"/>
                                  </w14:solidFill>
                                  <w14:prstDash w14:val="solid"/>
                                  <w14:round/>
                                </w14:textOutline>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:b/>
                                <w:color w:val="F7CAAC" w:themeColor="accent2" w:themeTint="66"/>
                                <w:sz w:val="72"/>
                                <w:szCs w:val="72"/>
                                <w14:textOutline w14:w="11112" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
                                  <w14:solidFill>
                                    <w14:schemeClr w14:val="accent2"/>
                                  </w14:solidFill>
                                  <w14:prstDash w14:val="solid"/>
                                  <w14:round/>
                                </w14:textOutline>
                              </w:rPr>
                              <w:t xml:space="preserve">Hola a todos trabajaremos con los </w:t>
                            </w:r>
                            <w:proofErr w:type="spellStart"/>
                            <w:r>
                              <w:rPr>
                                <w:b/>
                                <w:color w:val="F7CAAC" w:themeColor="accent2" w:themeTint="66"/>
                                <w:sz w:val="72"/>
                                <w:szCs w:val="72"/>
                                <w14:textOutline w14:w="11112" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
                                  <w14:solidFill>
                                    <w14:schemeClr w14:val="accent2"/>
                                  </w14:solidFill>
                                  <w14:prstDash w14:val="solid"/>
                                  <w14:round/>
                                </w14:textOutline>
                              </w:rPr>
                              <w:t>branch</w:t>
                            </w:r>
                            <w:proofErr w:type="spellEnd"/>
                          </w:p>
                        </w:txbxContent>
                      </wps:txbx>
                      <wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="none" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1">
                        <a:prstTxWarp prst="textNoShape">
                          <a:avLst/>
                        </a:prstTxWarp>
                        <a:spAutoFit/>
                      </wps:bodyPr>
                    </wps:wsp>
                  </a:graphicData>
                </a:graphic>
              </wp:anchor>
            </w:drawing>
          </mc:Choice>
          <mc:Fallback>
            <w:pict>
              <v:shapetype w14:anchorId="58297744" id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe">
                <v:stroke joinstyle="miter"/>
                <v:path gradientshapeok="t" o:connecttype="rect"/>
              </v:shapetype>
              <v:shape id="Cuadro de texto 2" o:spid="_x0000_s1026" type="#_x0000_t202" style="position:absolute;margin-left:0;margin-top:0;width:2in;height:2in;z-index:251659264;visibility:visible;mso-wrap-style:none;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQA/Ujy/JAIAAE4EAAAOAAAAZHJzL2Uyb0RvYy54bWysVE2P2jAQvVfqf7B8L4GItjQirCgrqkpo&#10;dyW22rNxbBLJ9li2IaG/vmMnYem2p6oXZ748nnnzJsu7TityFs43YEo6m0wpEYZD1ZhjSX88bz8s&#10;KPGBmYopMKKkF+Hp3er9u2VrC5FDDaoSjmAS44vWlrQOwRZZ5nktNPMTsMKgU4LTLKDqjlnlWIvZ&#10;tcry6fRT1oKrrAMuvEfrfe+kq5RfSsHDo5ReBKJKirWFdLp0HuKZrZasODpm64YPZbB/qEKzxuCj&#10;11T3LDBycs0fqXTDHXiQYcJBZyBlw0XqAbuZTd90s6+ZFakXBMfbK0z+/6XlD+cnR5qqpDklhmkc&#10;0ebEKgekEiSILgDJI0it9QXG7i1Gh+4rdDjs0e7RGHvvpNPxi10R9CPclyvEmInweGmRLxZTdHH0&#10;jQrmz16vW+fDNwGaRKGkDmeYoGXnnQ996BgSXzOwbZRKc1TmNwPmjJYs1t7XGKXQHbqhoQNUF+zH&#10;QU8Lb/m2wTd3zIcn5pAHWCdyOzziIRW0JYVBoqQG9/Nv9hiP40EvJS3yqqQGiU+J+m5wbF9m83mk&#10;YVLmHz/nqLhbz+HWY056A0jcGe6Q5UmM8UGNonSgX3AB1vFNdDHD8eWShlHchJ7ruEBcrNcpCIln&#10;WdiZveUxdYQs4vncvTBnB9Dj5B9g5B8r3mDfx8ab3q5PASeQBhPh7TEdUEfSptEOCxa34lZPUa+/&#10;gdUvAAAA//8DAFBLAwQUAAYACAAAACEAS4kmzdYAAAAFAQAADwAAAGRycy9kb3ducmV2LnhtbEyP&#10;0U7DMAxF35H4h8hIvLF0FaBSmk5owDMw+ACvMU1p41RNthW+HoOQxovlq2tdn1utZj+oPU2xC2xg&#10;uchAETfBdtwaeHt9vChAxYRscQhMBj4pwqo+PamwtOHAL7TfpFZJCMcSDbiUxlLr2DjyGBdhJBbv&#10;PUwek8ip1XbCg4T7QedZdq09diwfHI60dtT0m503UGT+qe9v8ufoL7+WV259Hx7GD2POz+a7W1CJ&#10;5nQ8hh98QYdamLZhxzaqwYAUSb9TvLwoRG7/Fl1X+j99/Q0AAP//AwBQSwECLQAUAAYACAAAACEA&#10;toM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQA&#10;BgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAAAAAAAC8BAABfcmVscy8ucmVsc1BLAQItABQA&#10;BgAIAAAAIQA/Ujy/JAIAAE4EAAAOAAAAAAAAAAAAAAAAAC4CAABkcnMvZTJvRG9jLnhtbFBLAQIt&#10;ABQABgAIAAAAIQBLiSbN1gAAAAUBAAAPAAAAAAAAAAAAAAAAAH4EAABkcnMvZG93bnJldi54bWxQ&#10;SwUGAAAAAAQABADzAAAAgQUAAAAA&#10;" filled="f" stroked="f">
                <v:fill o:detectmouseclick="t"/>
                <v:textbox style="mso-fit-shape-to-text:t">
                  <w:txbxContent>
                    <w:p>
                      <w:pPr>
                        <w:jc w:val="center"/>
                        <w:rPr>
                          <w:b/>
                          <w:color w:val="F7CAAC" w:themeColor="accent2" w:themeTint="66"/>
                          <w:sz w:val="72"/>
                          <w:szCs w:val="72"/>
                          <w14:textOutline w14:w="11112" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
                            <w14:solidFill>
                              <w14:schemeClr w14:val="accent2"/>
                            </w14:solidFill>
                            <w14:prstDash w14:val="solid"/>
                            <w14:round/>
                          </w14:textOutline>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:b/>
                          <w:color w:val="F7CAAC" w:themeColor="accent2" w:themeTint="66"/>
                          <w:sz w:val="72"/>
                          <w:szCs w:val="72"/>
                          <w14:textOutline w14:w="11112" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
                            <w14:solidFill>
                              <w14:schemeClr w14:val="accent2"/>
                            </w14:solidFill>
                            <w14:prstDash w14:val="solid"/>
                            <w14:round/>
                          </w14:textOutline>
                        </w:rPr>
                        <w:t xml:space="preserve">Hola a todos trabajaremos con los </w:t>
                      </w:r>
                      <w:proofErr w:type="spellStart"/>
                      <w:r>
                        <w:rPr>
                          <w:b/>
                          <w:color w:val="F7CAAC" w:themeColor="accent2" w:themeTint="66"/>
                          <w:sz w:val="72"/>
                          <w:szCs w:val="72"/>
                          <w14:textOutline w14:w="11112" w14:cap="flat" w14:cmpd="sng" w14:algn="ctr">
                            <w14:solidFill>
                              <w14:schemeClr w14:val="accent2"/>
                            </w14:solidFill>
                            <w14:prstDash w14:val="solid"/>
                            <w14:round/>
                          </w14:textOutline>
                        </w:rPr>
                        <w:t>branch</w:t>
                      </w:r>
                      <w:proofErr w:type="spellEnd"/>
                    </w:p>
                  </w:txbxContent>
                </v:textbox>
              </v:shape>
            </w:pict>
          </mc:Fallback>
        </mc:AlternateContent>
      </w:r>
    </w:p>
    </w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)
